$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, shifting existing rows 31-122 down to 32-123.
$ws.Rows(31).Insert()

# Populate the newly inserted row 31 with the new weekly data point
# (same dimension/variety as neighbouring rows, new date & price figures).
$ws.Cells.Item(31, 1).Value = 7
$ws.Cells.Item(31, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(31, 3).Value = "Ñuble"
$ws.Cells.Item(31, 4).Value = (Get-Date -Year 2023 -Month 3 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(31, 5).Value = 16
$ws.Cells.Item(31, 6).Value = 100112037
$ws.Cells.Item(31, 7).Value = "Cebollín"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 80
$ws.Cells.Item(31, 11).Value = 7000
$ws.Cells.Item(31, 12).Value = 7000
$ws.Cells.Item(31, 13).Value = 7000
$ws.Cells.Item(31, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(31, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(31, 16).Value = 194
$ws.Cells.Item(31, 17).Value = 36
$ws.Cells.Item(31, 18).Value = "Hortaliza"
